$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-7 got reshuffled (a cyclic permutation) while keeping
# column A,B,C,E,F,G,H,I,J,K fixed. Column L only changes for rows 6/7
# (they swap), for rows 2-5 it stays "Primera".
# New content per row for columns D,L,M,N,O,P,Q,R,S,T:

$rows = @{
    2 = @{ D = 44719; L = "Primera"; M = 50;  N = 14000; O = 15000; P = 14400; Q = "`$/caja 18 kilos granel";      R = "Región del Maule";      S = 800;   T = 18 }
    3 = @{ D = 44708; L = "Primera"; M = 70;  N = 12000; O = 13000; P = 12571; Q = "`$/caja 12 kilos empedrada";   R = "Provincia de Curicó";   S = 1048;  T = 12 }
    4 = @{ D = 44707; L = "Primera"; M = 60;  N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada";   R = "Provincia de Curicó";   S = 1042;  T = 12 }
    5 = @{ D = 44330; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel";      R = "Provincia de Curicó";   S = 861;   T = 18 }
    6 = @{ D = 44334; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/caja 12 kilos granel";      R = "Región de O'Higgins";   S = 11500; T = 1  }
    7 = @{ D = 44742; L = "Segunda"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos granel";      R = "Región de O'Higgins";   S = 806;   T = 18 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
